$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old single-row header data first
$ws.Cells.Clear()

# TABLE 1 title
$ws.Range("A1").Value = "TABLE 1"

# TABLE 1 column headers (row 2)
$ws.Range("A2").Value = "job_id"
$ws.Range("B2").Value = "job_type"
$ws.Range("C2").Value = "address"
$ws.Range("D2").Value = "customer_name"
$ws.Range("E2").Value = "job_status"
$ws.Range("F2").Value = "scheduled_date"
$ws.Range("G2").Value = "completed_on"
$ws.Range("H2").Value = "revenue"
$ws.Range("I2").Value = "on_site_hours"

# TABLE 2 title
$ws.Range("A4").Value = "TABLE 2: job_clock_events"

# TABLE 2 column headers (row 5)
$ws.Range("A5").Value = "id"
$ws.Range("B5").Value = "job_id"
$ws.Range("C5").Value = "tech_name"
$ws.Range("D5").Value = "hours"

# Selection as recorded in the target workbook
$ws.Range("B10:F14").Select()
